$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.00409836065573771
$ws.Cells.Item(2, 3).Value = 0.00307377049180328
$ws.Cells.Item(2, 4).Value = 0.00409836065573771
$ws.Cells.Item(2, 5).Value = 0.00717213114754098
$ws.Cells.Item(2, 6).Value = 0.00102459016393443
$ws.Cells.Item(2, 7).Value = 0.00409836065573771
$ws.Cells.Item(2, 8).Value = 0.00307377049180328
$ws.Cells.Item(2, 9).Value = 0.839139344262295
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0.00102459016393443
$ws.Cells.Item(2, 12).Value = 0.00102459016393443
$ws.Cells.Item(2, 13).Value = 0.00409836065573771
$ws.Cells.Item(2, 14).Value = 0.00102459016393443
$ws.Cells.Item(2, 15).Value = 0.00102459016393443
$ws.Cells.Item(2, 16).Value = 0.00717213114754098
$ws.Cells.Item(2, 17).Value = 0.211065573770492
$ws.Cells.Item(2, 18).Value = 0.00307377049180328
$ws.Cells.Item(2, 19).Value = 0.98155737704918
$ws.Cells.Item(2, 20).Value = 0.983606557377049
$ws.Cells.Item(2, 21).Value = 0.040983606557377
$ws.Cells.Item(2, 22).Value = 0.00409836065573771
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0.978483606557377

$ws.Cells.Item(3, 2).Value = 0.862704918032787
$ws.Cells.Item(3, 3).Value = 0.939549180327869
$ws.Cells.Item(3, 4).Value = 0.0153688524590164
$ws.Cells.Item(3, 5).Value = 0.00204918032786885
$ws.Cells.Item(3, 6).Value = 0.0245901639344262
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0.040983606557377
$ws.Cells.Item(3, 10).Value = 0.00204918032786885
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0.00307377049180328
$ws.Cells.Item(3, 13).Value = 0.00307377049180328
$ws.Cells.Item(3, 14).Value = 0.00307377049180328
$ws.Cells.Item(3, 15).Value = 0.986680327868853
$ws.Cells.Item(3, 16).Value = 0.00307377049180328
$ws.Cells.Item(3, 17).Value = 0.00512295081967213
$ws.Cells.Item(3, 18).Value = 0.973360655737705
$ws.Cells.Item(3, 19).Value = 0.00819672131147541
$ws.Cells.Item(3, 20).Value = 0.00614754098360656
$ws.Cells.Item(3, 21).Value = 0.00204918032786885
$ws.Cells.Item(3, 22).Value = 0.857581967213115
$ws.Cells.Item(3, 23).Value = 0.00512295081967213
$ws.Cells.Item(3, 24).Value = 0.00204918032786885

$ws.Cells.Item(4, 2).Value = 0.125
$ws.Cells.Item(4, 3).Value = 0.00307377049180328
$ws.Cells.Item(4, 4).Value = 0.00102459016393443
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0.00717213114754098
$ws.Cells.Item(4, 7).Value = 0.995901639344262
$ws.Cells.Item(4, 8).Value = 0.994877049180328
$ws.Cells.Item(4, 9).Value = 0.111680327868852
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0.997950819672131
$ws.Cells.Item(4, 12).Value = 0.995901639344262
$ws.Cells.Item(4, 13).Value = 0.00409836065573771
$ws.Cells.Item(4, 14).Value = 0.992827868852459
$ws.Cells.Item(4, 15).Value = 0.00307377049180328
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0.00409836065573771
$ws.Cells.Item(4, 18).Value = 0.00409836065573771
$ws.Cells.Item(4, 19).Value = 0.00307377049180328
$ws.Cells.Item(4, 20).Value = 0.00204918032786885
$ws.Cells.Item(4, 21).Value = 0.945696721311475
$ws.Cells.Item(4, 22).Value = 0.00307377049180328
$ws.Cells.Item(4, 23).Value = 0.00717213114754098
$ws.Cells.Item(4, 24).Value = 0.0163934426229508

$ws.Cells.Item(5, 2).Value = 0.00717213114754098
$ws.Cells.Item(5, 3).Value = 0.0543032786885246
$ws.Cells.Item(5, 4).Value = 0.979508196721312
$ws.Cells.Item(5, 5).Value = 0.99077868852459
$ws.Cells.Item(5, 6).Value = 0.967213114754098
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0.00102459016393443
$ws.Cells.Item(5, 9).Value = 0.00819672131147541
$ws.Cells.Item(5, 10).Value = 0.997950819672131
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0.988729508196721
$ws.Cells.Item(5, 14).Value = 0.00307377049180328
$ws.Cells.Item(5, 15).Value = 0.00922131147540984
$ws.Cells.Item(5, 16).Value = 0.989754098360656
$ws.Cells.Item(5, 17).Value = 0.779713114754098
$ws.Cells.Item(5, 18).Value = 0.0194672131147541
$ws.Cells.Item(5, 19).Value = 0.00717213114754098
$ws.Cells.Item(5, 20).Value = 0.00819672131147541
$ws.Cells.Item(5, 21).Value = 0.0102459016393443
$ws.Cells.Item(5, 22).Value = 0.135245901639344
$ws.Cells.Item(5, 23).Value = 0.987704918032787
$ws.Cells.Item(5, 24).Value = 0.00307377049180328
